$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WAT url")

# --- Capture the text currently sitting in A11 (it is about to be moved to K2) ---
$movedText = $ws.Cells(11,1).Text

# --- Remove all existing hyperlink objects; we will re-create them below, pointed
#     at their post-shift cell addresses (this engine does not auto-shift hyperlink
#     ranges when rows are deleted) ---
$ws.Hyperlinks.Delete()

# --- Relocate the "promotion" URL text that was sitting in A11 over to K2 ---
$ws.Range("K2").Value = $movedText
$ws.Range("K2").Style = $ws.Range("A11").Style

# --- Delete row 11 outright; rows 12-37 shift up to become rows 11-36 ---
$ws.Rows("11:11").Delete()

# --- Re-add the 21 original hyperlinks, now pointing at their shifted cells ---
$ws.Hyperlinks.Add($ws.Cells(2,1), "https://www.watsons.com.hk/neutrogena-deep-clean-brightening-foaming-cleanser-100g/p/BP_261330")
$ws.Hyperlinks.Add($ws.Cells(4,1), "https://www.watsons.com.hk/cleansing-sheet-moist/p/BP_250613")
$ws.Hyperlinks.Add($ws.Cells(3,1), "https://www.watsons.com.hk/bifesta-cleansing-lotion-moist/p/BP_296288")
$ws.Hyperlinks.Add($ws.Cells(7,1), "https://www.watsons.com.hk/perfect-white-clay-120g/p/BP_292971")
$ws.Hyperlinks.Add($ws.Cells(9,1), "https://www.watsons.com.hk/sensibio-h2o-watsons/p/BP_262653")
$ws.Hyperlinks.Add($ws.Cells(10,1), "https://www.watsons.com.hk/bioderma-sensibio-bonus-set-pump/p/BP_807961")
$ws.Hyperlinks.Add($ws.Cells(16,1), "https://www.watsons.com.hk/sunplay-super-block-spf130pa/p/BP_801308")
$ws.Hyperlinks.Add($ws.Cells(15,1), "https://www.watsons.com.hk/sunplay-skin-aqua-phyiscal-sunscreen-spf50/p/BP_801337")
$ws.Hyperlinks.Add($ws.Cells(18,1), "https://www.watsons.com.hk/ultra-mist-sport-sunscreen-lotion-spf110-90ml/p/BP_245822")
$ws.Hyperlinks.Add($ws.Cells(17,1), "https://www.watsons.com.hk/aloe-vera-gel-230g/p/BP_214152")
$ws.Hyperlinks.Add($ws.Cells(26,1), "https://www.watsons.com.hk/ptn-miracles-crystal-smooth-trt500g/p/BP_370593")
$ws.Hyperlinks.Add($ws.Cells(25,1), "https://www.watsons.com.hk/rerise-hair-colour-treatment-reblack-tame/p/BP_261543")
$ws.Hyperlinks.Add($ws.Cells(27,1), "https://www.watsons.com.hk/lux-hair-supplement-shampoo-450g-natural-shine/p/BP_805400")
$ws.Hyperlinks.Add($ws.Cells(28,1), "https://www.watsons.com.hk/micellar-volume-sh/p/BP_406820")
$ws.Hyperlinks.Add($ws.Cells(29,1), "https://www.watsons.com.hk/moist-diane-miracle-you-shampoo/p/BP_801668")
$ws.Hyperlinks.Add($ws.Cells(30,1), "https://www.watsons.com.hk/moist-diane-perfect-beauty-extra-hair-fall-control-shampoo-450ml/p/BP_806076")
$ws.Hyperlinks.Add($ws.Cells(31,1), "https://www.watsons.com.hk/color-care-shampoo/p/BP_296291")
$ws.Hyperlinks.Add($ws.Cells(32,1), "https://www.watsons.com.hk/50-megumi-anti-grey-shampoo-400ml/p/BP_805721")
$ws.Hyperlinks.Add($ws.Cells(33,1), "https://www.watsons.com.hk/damage-care-shampoo/p/BP_242363")
$ws.Hyperlinks.Add($ws.Cells(35,1), "https://www.watsons.com.hk/hair-recipe-kiwi-fig-volume-shampoo-530ml/p/BP_807103")
$ws.Hyperlinks.Add($ws.Cells(36,1), "https://www.watsons.com.hk/moisturizing-ad-sh/p/BP_244823")

# --- New hyperlink (rId22): A11 now holds the former A12 text/value; it previously
#     had no hyperlink of its own, so "Insert Hyperlink" has been applied to it ---
$ws.Hyperlinks.Add($ws.Cells(11,1), "https://www.watsons.com.hk/perfect-uv-sunscreen-skincare-gel-spf50-pa-90g/p/BP_801291")

# Restore the plain (non-highlighted) Hyperlink cell style on the newly-linked A11
# cell -- Hyperlinks.Add() otherwise stamps a synthetic style variant on it.
$ws.Cells(11,1).Style = "Hyperlink"

# --- Selection / view bookkeeping to match the post-edit window state ---
$ws.Range("L9").Select()
